# Apply updates to Alpha2F-HW20 worksheet:
#  - insert new shared strings / labels (OffsetF, OffsetA, RD Single, TD Single, 1Pair-B)
#  - extend grid from 17 data columns (C:S) to 18 data columns (C:T)
#  - recompute/replace the numeric data for rows 3-7 (Equal Angle..OffsetF)
#  - append 4 new rows (OffsetA, RD Single, TD Single, HexGrid-90degTilt5degRes)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, extend the header/style formatting to the new column (T) and new rows (8-11),
# matching the existing bold/bordered style used for the header row and the row labels
# in column A.
$ws.Range("S1").Copy($ws.Range("T1"))
$ws.Range("A2").Copy($ws.Range("A8:A11"))

# Full target grid (row 1 = A1:T1, row 2 = A2:T2, ... row 11 = A11:T11)
$data = @(
    @($null, 0, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18),
    @(0, "HKL", "[1, 1, 0]", "[2, 0, 0]", "[2, 1, 1]", "[2, 2, 0]", "[3, 1, 0]", "[2, 2, 2]", "[3, 2, 1]", "[4, 0, 0]", "1Pair-A", "1Pair-B", "2Pairs-A", "2Pairs-B", "3Pairs-A", "3Pairs-B", "3Pairs-C", "4Pairs", "5A4F", "MaxUnique"),
    @(1, "Equal Angle", 0.959935158501441, 1.049798270893372, 0.9871181556195966, 0.959935158501441, 1.036909221902017, 0.9421037463976946, 0.979985590778098, 1.049798270893372, 0.959935158501441, 0.9871181556195966, 1.018458213256484, 1.018458213256484, 1.024608549471662, 0.9989505283381365, 0.9989505283381365, 0.9891966858789626, 0.9891966858789626, 0.9926416906820363),
    @(2, "CLR", 0.9962427724191131, 0.9962065869928907, 0.9933546920473233, 0.9962427724191131, 0.9991383552946341, 0.9876755546390646, 0.9939013271023057, 0.9962065869928907, 0.9962427724191131, 0.9933546920473233, 0.994780639520107, 0.994780639520107, 0.9962332114449494, 0.995268017153109, 0.9952680171531091, 0.9955117059696101, 0.9955117059696101, 0.9944198814158886),
    @(3, "BT8Hex", 0.9913555945842842, 1.002565099727697, 0.9924645218287246, 0.9913555945842842, 1.003847983899797, 0.9824777824418984, 0.9918234313627076, 1.002565099727697, 0.9913555945842842, 0.9924645218287246, 0.9975148107782108, 0.9975148107782108, 0.9996258684854062, 0.9954617387135686, 0.9954617387135687, 0.9944352026812475, 0.9944352026812475, 0.9940890689741848),
    @(4, "Spiral", 0.9960959314838964, 0.9944566412605269, 0.994848783720965, 0.9960959314838964, 0.9962912928802704, 0.9936107770493551, 0.9949390778223987, 0.9944566412605269, 0.9960959314838964, 0.994848783720965, 0.9946527124907459, 0.9946527124907459, 0.9951989059539207, 0.9951337854884628, 0.9951337854884628, 0.9953743219873212, 0.9953743219873212, 0.9950404173695687),
    @(5, "OffsetF", 1.387339174743095, 0.39716511093343, 1.024481114286854, 1.387339174743095, 0.7111694987268753, 1.066063385505465, 1.129411962922567, 0.39716511093343, 1.387339174743095, 1.024481114286854, 0.7108231126101417, 0.7108231126101417, 0.710938574649053, 0.9363284666544596, 0.9363284666544596, 1.049081143676618, 1.049081143676618, 0.952605041186381),
    @(6, "OffsetA", 0.8958226795379629, 0.8737749544024488, 1.079422537393283, 0.8958226795379629, 0.9352340750914031, 1.180947956254795, 1.026525138625567, 0.8737749544024488, 0.8958226795379629, 1.079422537393283, 0.976598745897866, 0.976598745897866, 0.9628105222957117, 0.949673390444565, 0.949673390444565, 0.9362107127179145, 0.9362107127179145, 0.9986212235509101),
    @(7, "RD Single", 1.97, 0.21, 0.83, 1.97, 0.63, 0.6899999999999999, 1.14, 0.21, 1.97, 0.83, 0.52, 0.52, 0.5566666666666666, 1.003333333333333, 1.003333333333333, 1.245, 1.245, 0.9116666666666666),
    @(8, "TD Single", 0.84, 0.07000000000000001, 1.42, 0.84, 0.29, 2.44, 1.27, 0.07000000000000001, 0.84, 1.42, 0.745, 0.745, 0.5933333333333334, 0.7766666666666667, 0.7766666666666667, 0.7925, 0.7925, 1.055),
    @(9, "HexGrid-90degTilt5degRes", 0.9977772906291448, 0.99357208486699, 0.9944692573089144, 0.9977772906291448, 0.9959429918347955, 0.993069502023931, 0.9950548353801498, 0.99357208486699, 0.9977772906291448, 0.9944692573089144, 0.9940206710879522, 0.9940206710879522, 0.9946614446702333, 0.9952728776016831, 0.9952728776016831, 0.9958989808585486, 0.9958989808585486, 0.9949809936739875)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $v = $row[$c]
        if ($null -ne $v) {
            $ws.Cells.Item($r + 1, $c + 1).Value = $v
        }
    }
}

Write-Host "Alpha2F-HW20 sheet updated: extended to A1:T11 with OffsetF/OffsetA/RD Single/TD Single/1Pair-B rows+cols."
